# furnace_age_EIA.xlsx -- "updated figures 5 and 6"
#
# Relabels the purchase-year buckets used in figures 5/6 and widens the
# 2007-2009 bucket's upper formula weight, then tidies up a block of
# left-over blank rows below the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the purchase-year bucket (both the 2009-RECS and 2015-RECS
#     copies of the label) from "(2007-2009]" / "[2008-2010)" to "[2008-2009]" ---
$ws.Range("A7").Value = "[2008-2009]"
$ws.Range("A12").Value = "[2008-2009]"

# --- Rows 15:16 and 20:27 are blank spacer rows left over under the table;
#     rows 14, 17, 18 and 19 hold the "2000 and earlier" / summary figures.
#     Delete the whole 14:27 block (clears the stray formatting on it) and
#     rebuild the rows that actually carry data so 14/17/18/19 land back on
#     the same row numbers. ---
$ws.Rows("14:27").Delete()

# Row 14 - "2000 and earlier" (Main Heating), continuing the B13:B14 pattern
$ws.Range("A14").Value = "2000 and earlier"
$ws.Range("B14").Formula = '=B9*$B$4/$B$3'

# Row 17 - "Main NG Central Warm-Air Furnaces" sub-heading
$ws.Range("A17").Value = "Main NG Central Warm-Air Furnaces"

# Row 18 - "Purchased before 2007" summary figure
$ws.Range("A18").Value = "Purchased before 2007"
$ws.Range("B18").Formula = '=B13*2/3+B14'

# Row 19 - relabelled summary figure, now dividing the 2015 RECS term by 5
# instead of 4, with its wrapped two-line row height restored
$ws.Range("A19").Value = "Purchased after and including 2007, before 2010"
$ws.Range("A19").WrapText = $true
$ws.Rows(19).RowHeight = 51
$ws.Range("B19").Formula = '=B13/3+B12+E12/5'

# --- Active cell/selection moved to H11 ---
$ws.Range("H11").Select() | Out-Null
